$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Rows 18-20: new year rows (Ano + 3 numeric columns), matching the
# look/format (style) already used by the existing data rows (e.g. row 17).
# ---------------------------------------------------------------------

# Copy the formatting of A17 (style index 1: centered "General" xf) onto
# the full A18:D20 block so every new cell - including column B, which in
# the existing rows uses a different currency-ish style - lands on that
# same style.
$ws.Range("A17").Copy()
$ws.Range("A18:D20").PasteSpecial(-4122)

# Column A must hold the year as TEXT (shared string), not a number, so
# force text formatting before typing the value, then restore the cell's
# normal/general format/style afterwards (without touching the entered text).
$ws.Range("A18:A20").NumberFormat = "@"
$ws.Range("A18").Value = "1971"
$ws.Range("A19").Value = "1972"
$ws.Range("A20").Value = "1973"
$ws.Range("A17").Copy()
$ws.Range("A18:A20").PasteSpecial(-4122)

# Numeric data for the three value columns.
$ws.Range("B18").Value = 4.92
$ws.Range("C18").Value = 5.15
$ws.Range("D18").Value = 5.08

$ws.Range("B19").Value = 4.92
$ws.Range("C19").Value = 5.15
$ws.Range("D19").Value = 5.08

$ws.Range("B20").Value = 1.63
$ws.Range("C20").Value = 1.97
$ws.Range("D20").Value = 1.02

# ---------------------------------------------------------------------
# Rows 21-29: blank cells (B:D only) carrying the same style as the data
# rows above them, mirroring a fill/format-only extension of the table.
# ---------------------------------------------------------------------
$ws.Range("A17").Copy()
$ws.Range("B21:D29").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Extend the sheet's recognised/used extent down to row 1000 (touch a
# single, still-empty, General-formatted cell so nothing visible changes).
# ---------------------------------------------------------------------
$ws.Range("D1000").NumberFormat = "General"

# Final selection, matching the saved view state.
$ws.Range("D21").Select()
